# "Team Members Had Space" - remove the stray leading space/nbsp before
# "Carla Machado" on the Team Members slide, re-splitting the bold name
# run into "Carla " + "Machado " (matching how PowerPoint re-flows the
# run boundaries once the leading space character is deleted).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Remove the leading space (non-breaking space) that preceded "Carla Machado".
$leading = $tr.Characters(1,1)
$leading.Text = ""

# Re-assert the bold formatting on "Carla " (first 6 characters of the now
# merged "Carla Machado " run) which causes it to split into its own run,
# matching the target "Carla " / "Machado " run split.
$carla = $tr.Characters(1,6)
$carla.Font.Bold = -1
